$d = $word.ActiveDocument

# 1. Remove the old hidden "_GoBack" bookmark that currently sits right
#    after "...quickly." at the end of the file-IO paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the end of the "Reflecting on..." paragraph (the prompt text
#    under the "Personal Evaluation" heading) and insert the new
#    reflection content directly after it, before the "References"
#    heading. The inserted XML fragment also re-creates the "_GoBack"
#    bookmark in its new location (inside the final "Overall," paragraph).
$findRng = $d.Content
$found = $findRng.Find.Execute("Reflecting on what you have learned, the challenges you faced, the methods you used to overcome challenges and how you feel you have performed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Reflecting on...' paragraph"
}
$insertPos = $findRng.End
$insertRng = $d.Range($insertPos, $insertPos)

$newContentXml = '    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">What </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>ive</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> leaned</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>arrays</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">The first thing I learned to complete this coursework is how arrays work. I already had a basic understanding of this but for this coursework I needed to learn ways to display and add to them in ways that I never had before. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>An example of this is my insert function, it inserts a value at a position in an array, this taught me how to add to an array in a concise way by calling a function rather than creating a loop to add to an array each time I needed to do so.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Stacks/structs</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">The second type of array I needed to learn is a two </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>dimentional</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> array, I used this in my stack. This took some time to get my head round as I had never used one in this way before. For my stack I needed to use the two </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>dimentional</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> array to store my list of boards and to add and remove them </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>fronm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> the list I had to re-think how I worked with arrays. Making sure to keep the column </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>consistant</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> while I looped through each element of the row. The lab also helped me to learn about stacks as well, the push, initialize and pop functions were all easy to understand once I had got my head around how id implement them in a two </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>dimentional</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> array rather than a one </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>dimentional</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> array. The last part of my stack I had to learn was a struct, I had </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>come into contact with</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> this datat</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>y</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>pe before but using it</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> to create a stack. </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>However</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> using it as a stack helped me to use it just like any other data type, passing it around as a parameter into other </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>funcitons</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Mastered functions</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">The final thing I improved is my sills with functions. This </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>cousework</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> helped me to learn a lot about them and features such as passing </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>peramaters</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> and</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> changing values that are declared in the main within a function. The passing of parameters through a function was the first thing I needed to learn as the easiest place to start this coursework was having a display function working </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">and to get this working I needed to pass an array to be printed. The next thing I needed a function to do was change a value that is declared in the main. This was done in the check winner function as the winner integer is the condition of my main while loop in the main function, and </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>in order to</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> break that loop I needed to change the value of the integer from zero to one.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Challenges</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Stack</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">The biggest challenge I faced was implementing the stack. I needed to learn both how </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>my</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> struct worked and how that related to my two </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>dimentional</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> array contained within it. After studying the stack implementation in the </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>lab</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> I eventually realized that the top integer declared in my struct should be used to keep track of the row in my array, this would be the number of each ‘layer’ of the stack. This helped me to implement both the push and pop functions as I realized that all I needed to do was loop through the columns of each row and only increment or decrement the top value once per function. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Breaking the loop when there’s a win/draw</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">The next challenge I faced was making sure in the main method the user was only asked for moves when it was valid, i.e. there is no winner and the board isn’t full. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">To overcome </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>this</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> I used two variables to make sure that if there was a winner or the board became full the while loop would break and the user wouldn’t be asked for any more moves. These variables get passed to functions which check for both of these events, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>checkWinner</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">) and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>checkFullBoard</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>() and when either one is true the while loop breaks and a set of if-else statements declare what will happen afterwards.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">Performed well overall, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>fulfilled</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> most of the requirements</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Overall</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> I think I performed well, I reached all of the minimum requirements set out in the description </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Arial"/>
          <w:sz w:val="20"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>and I went above them when I included my undo feature and my replaying an old game feature. Alongside this my program runs quickly and efficiently and is easy to understand and use for a new user.</w:t>
      </w:r>
    </w:p>'

$insertRng.InsertXML($newContentXml)
